$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'42.973.08"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.74%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.540.92"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.29%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.05%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'318.77"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.98%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'97.61"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +2.38%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.575"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -0.96%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D9').Value = "'0.537"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -0.63%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'36.44"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +0.40%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.65%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'7.63"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -1.44%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = "'  -3.12%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'2.933.88"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +0.55%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'2.531.28"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -0.96%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'15.19"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -3.26%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'  -0.84%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'43.051.29"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +0.84%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'6.86"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +3.20%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'12.82"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -1.92%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'0.0₃0972"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +0.10%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'70.00"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -1.90%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'254.17"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.13%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = "'  -0.32%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'2.06"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +0.96%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'26.57"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -3.79%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  +0.52%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  +4.52%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'40.74"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +2.92%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'10.47"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +3.97%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'5.94"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +0.11%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'158.43"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +1.19%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'2.17"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +2.35%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  +0.23%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  +3.92%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'19.07"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -4.58%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'0.0794"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +0.77%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  -0.17%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'2.49"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +12.75%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  -0.75%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'22.15"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -9.94%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'3.85"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -0.06%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'0.0304"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.04%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = "'  +0.31%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'3.30"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -3.06%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'2.019.15"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -1.55%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'9.14"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +2.81%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'84.58"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -1.63%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'76.62"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +2.70%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'106.56"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +4.21%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'2.790.49"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +0.74%  "
$ws.Range('E51').Style = 'Normal'
